# Add two new time-log rows (row 3 and row 4) to Sheet1, matching the
# style of the existing date row, widen column D to fit the new text,
# and move the active selection the way the author's session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 -----------------------------------------------------------
# Copy A2's format (date number format m/d/yyyy => style index 1) onto
# A3:A4 before writing the values, so the new date cells share the same
# style as the existing one instead of Excel minting a brand-new style.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)

$ws.Range("A3").Value = 43865
$ws.Range("B3").Value = "jaclemon"
$ws.Range("C3").Value = "60 minutes"
$ws.Range("D3").Value = "Wrote function to read command line arguments using notepad++"

# --- Row 4 -------------------------------------------------------------
$ws.Range("A4").Value = 43866
$ws.Range("B4").Value = "jaclemon"
$ws.Range("C4").Value = "30 minutes"
$ws.Range("D4").Value = "Used notepad++ to fix exit status and deleted branch and added all changes to master with github"

# --- Widen column D so the longer description text fits ----------------
$ws.Columns.Item(4).ColumnWidth = 54.5

# --- Move the selection to where the author left off --------------------
$ws.Range("D13").Select()
